# Commit message: "Fruta / hortaliza, semanal"
#
# A new weekly price record was inserted for
# "Terminal La Palmera de La Serena - Zanahoria" as row 284
# (pushing every subsequent row down by one, row 358 -> 359,
# and extending the used range from A1:R358 to A1:R359).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before the current row 284, shifting
# rows 284..358 down to 285..359.
$ws.Rows.Item(284).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(284, 1).Value  = 8
$ws.Cells.Item(284, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(284, 3).Value  = "Coquimbo"
$ws.Cells.Item(284, 4).Value  = 44722
$ws.Cells.Item(284, 5).Value  = 4
$ws.Cells.Item(284, 6).Value  = 100114013
$ws.Cells.Item(284, 7).Value  = "Zanahoria"
$ws.Cells.Item(284, 8).Value  = "Sin especificar"
$ws.Cells.Item(284, 9).Value  = "Primera"
$ws.Cells.Item(284, 10).Value = 600
$ws.Cells.Item(284, 11).Value = 6000
$ws.Cells.Item(284, 12).Value = 7000
$ws.Cells.Item(284, 13).Value = 6500
$ws.Cells.Item(284, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(284, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(284, 16).Value = 325
$ws.Cells.Item(284, 17).Value = 20
$ws.Cells.Item(284, 18).Value = "Hortaliza"
